# Updates cryptos list values (price/volume columns, plus a handful of
# coin-name/link/price/volume row swaps) per the upstream data refresh.
#
# Values that look numeric (e.g. "0.520", "18.28") must stay stored as TEXT
# (matching the workbook's existing inlineStr convention for column D/E) so
# trailing zeros and thousand-dot-separated figures like "60.508.95" are not
# mangled into real numbers. We force text storage by switching the cell to
# the "@" (Text) number format before writing the value, then clear the
# format again afterwards so no stray style/number-format is left behind on
# cells that were plain/default-styled before this edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param($Cell, [string]$Text)
    $Cell.NumberFormat = "@"
    $Cell.Value = $Text
    $Cell.ClearFormats()
}


# Row 2
Set-TextCell ($ws.Cells.Item(2, 4)) '60.508.95'
Set-TextCell ($ws.Cells.Item(2, 5)) '  -0.87%  '

# Row 3
Set-TextCell ($ws.Cells.Item(3, 4)) '2.350.36'
Set-TextCell ($ws.Cells.Item(3, 5)) '  -4.15%  '

# Row 4
Set-TextCell ($ws.Cells.Item(4, 5)) '  -0.12%  '

# Row 5
Set-TextCell ($ws.Cells.Item(5, 4)) '539.46'
Set-TextCell ($ws.Cells.Item(5, 5)) '  -0.92%  '

# Row 6
Set-TextCell ($ws.Cells.Item(6, 4)) '135.76'
Set-TextCell ($ws.Cells.Item(6, 5)) '  -6.27%  '

# Row 7
Set-TextCell ($ws.Cells.Item(7, 5)) '  -0.05%  '

# Row 8
Set-TextCell ($ws.Cells.Item(8, 4)) '0.520'
Set-TextCell ($ws.Cells.Item(8, 5)) '  -11.41%  '

# Row 9
Set-TextCell ($ws.Cells.Item(9, 4)) '2.349.18'
Set-TextCell ($ws.Cells.Item(9, 5)) '  -4.16%  '

# Row 10
Set-TextCell ($ws.Cells.Item(10, 4)) '0.104'
Set-TextCell ($ws.Cells.Item(10, 5)) '  -1.43%  '

# Row 11
Set-TextCell ($ws.Cells.Item(11, 5)) '  +0.16%  '

# Row 12
Set-TextCell ($ws.Cells.Item(12, 4)) '5.22'
Set-TextCell ($ws.Cells.Item(12, 5)) '  -3.62%  '

# Row 13
Set-TextCell ($ws.Cells.Item(13, 4)) '0.339'
Set-TextCell ($ws.Cells.Item(13, 5)) '  -3.09%  '

# Row 14
Set-TextCell ($ws.Cells.Item(14, 4)) '24.44'
Set-TextCell ($ws.Cells.Item(14, 5)) '  -5.36%  '

# Row 15
Set-TextCell ($ws.Cells.Item(15, 4)) '2.770.39'
Set-TextCell ($ws.Cells.Item(15, 5)) '  -4.38%  '

# Row 16
Set-TextCell ($ws.Cells.Item(16, 4)) '60.781.94'
Set-TextCell ($ws.Cells.Item(16, 5)) '  -0.36%  '

# Row 17
Set-TextCell ($ws.Cells.Item(17, 4)) '0.0000161'
Set-TextCell ($ws.Cells.Item(17, 5)) '  -2.50%  '

# Row 18
Set-TextCell ($ws.Cells.Item(18, 4)) '2.352.09'
Set-TextCell ($ws.Cells.Item(18, 5)) '  -4.28%  '

# Row 19
Set-TextCell ($ws.Cells.Item(19, 4)) '10.54'
Set-TextCell ($ws.Cells.Item(19, 5)) '  -4.45%  '

# Row 20
Set-TextCell ($ws.Cells.Item(20, 4)) '313.20'
Set-TextCell ($ws.Cells.Item(20, 5)) '  -1.46%  '

# Row 21
Set-TextCell ($ws.Cells.Item(21, 4)) '4.04'

# Row 22
Set-TextCell ($ws.Cells.Item(22, 4)) '6.56'
Set-TextCell ($ws.Cells.Item(22, 5)) '  -6.27%  '

# Row 23
Set-TextCell ($ws.Cells.Item(23, 4)) '5.82'
Set-TextCell ($ws.Cells.Item(23, 5)) '  -1.80%  '

# Row 24
Set-TextCell ($ws.Cells.Item(24, 5)) '  -0.09%  '

# Row 25
Set-TextCell ($ws.Cells.Item(25, 4)) '1.87'
Set-TextCell ($ws.Cells.Item(25, 5)) '  +1.39%  '

# Row 26
Set-TextCell ($ws.Cells.Item(26, 4)) '62.85'
Set-TextCell ($ws.Cells.Item(26, 5)) '  -1.16%  '

# Row 27
Set-TextCell ($ws.Cells.Item(27, 4)) '8.42'
Set-TextCell ($ws.Cells.Item(27, 5)) '  +9.85%  '

# Row 28
Set-TextCell ($ws.Cells.Item(28, 5)) '  +0.26%  '

# Row 29
Set-TextCell ($ws.Cells.Item(29, 4)) '2.464.17'
Set-TextCell ($ws.Cells.Item(29, 5)) '  -4.32%  '

# Row 30
Set-TextCell ($ws.Cells.Item(30, 2)) 'PEPE'
Set-TextCell ($ws.Cells.Item(30, 3)) 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
Set-TextCell ($ws.Cells.Item(30, 4)) '0.0₃0886'
Set-TextCell ($ws.Cells.Item(30, 5)) '  -7.72%  '

# Row 31
Set-TextCell ($ws.Cells.Item(31, 2)) 'InternetComputer(DFINITY)'
Set-TextCell ($ws.Cells.Item(31, 3)) 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextCell ($ws.Cells.Item(31, 4)) '7.87'
Set-TextCell ($ws.Cells.Item(31, 5)) '  -3.54%  '

# Row 32
Set-TextCell ($ws.Cells.Item(32, 2)) 'Fetch.AI'
Set-TextCell ($ws.Cells.Item(32, 3)) 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
Set-TextCell ($ws.Cells.Item(32, 4)) '1.37'
Set-TextCell ($ws.Cells.Item(32, 5)) '  -6.24%  '

# Row 33
Set-TextCell ($ws.Cells.Item(33, 2)) 'Bittensor'
Set-TextCell ($ws.Cells.Item(33, 3)) 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
Set-TextCell ($ws.Cells.Item(33, 4)) '496.48'
Set-TextCell ($ws.Cells.Item(33, 5)) '  -7.26%  '

# Row 34
Set-TextCell ($ws.Cells.Item(34, 5)) '  -2.37%  '

# Row 35
Set-TextCell ($ws.Cells.Item(35, 4)) '1.77'
Set-TextCell ($ws.Cells.Item(35, 5)) '  -5.82%  '

# Row 36
Set-TextCell ($ws.Cells.Item(36, 5)) '  -4.10%  '

# Row 37
Set-TextCell ($ws.Cells.Item(37, 4)) '0.999'
Set-TextCell ($ws.Cells.Item(37, 5)) '  -0.17%  '

# Row 38
Set-TextCell ($ws.Cells.Item(38, 4)) '4.53'
Set-TextCell ($ws.Cells.Item(38, 5)) '  -5.39%  '

# Row 39
Set-TextCell ($ws.Cells.Item(39, 4)) '0.370'
Set-TextCell ($ws.Cells.Item(39, 5)) '  -1.32%  '

# Row 40
Set-TextCell ($ws.Cells.Item(40, 4)) '18.28'
Set-TextCell ($ws.Cells.Item(40, 5)) '  -0.22%  '

# Row 41
Set-TextCell ($ws.Cells.Item(41, 4)) '5.18'
Set-TextCell ($ws.Cells.Item(41, 5)) '  -8.68%  '

# Row 42
Set-TextCell ($ws.Cells.Item(42, 4)) '1.76'
Set-TextCell ($ws.Cells.Item(42, 5)) '  +0.92%  '

# Row 43
Set-TextCell ($ws.Cells.Item(43, 2)) 'USDe'
Set-TextCell ($ws.Cells.Item(43, 3)) 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
Set-TextCell ($ws.Cells.Item(43, 4)) '1.00'
Set-TextCell ($ws.Cells.Item(43, 5)) '  +0.00%  '

# Row 44
Set-TextCell ($ws.Cells.Item(44, 2)) 'Monero'
Set-TextCell ($ws.Cells.Item(44, 3)) 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextCell ($ws.Cells.Item(44, 4)) '138.10'
Set-TextCell ($ws.Cells.Item(44, 5)) '  -3.60%  '

# Row 45
Set-TextCell ($ws.Cells.Item(45, 4)) '40.04'
Set-TextCell ($ws.Cells.Item(45, 5)) '  -1.03%  '

# Row 46
Set-TextCell ($ws.Cells.Item(46, 2)) 'Aave'
Set-TextCell ($ws.Cells.Item(46, 3)) 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextCell ($ws.Cells.Item(46, 4)) '140.14'
Set-TextCell ($ws.Cells.Item(46, 5)) '  -1.17%  '

# Row 47
Set-TextCell ($ws.Cells.Item(47, 2)) 'dogwifhat'
Set-TextCell ($ws.Cells.Item(47, 3)) 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
Set-TextCell ($ws.Cells.Item(47, 4)) '2.08'
Set-TextCell ($ws.Cells.Item(47, 5)) '  -7.61%  '

# Row 48
Set-TextCell ($ws.Cells.Item(48, 4)) '3.49'
Set-TextCell ($ws.Cells.Item(48, 5)) '  -2.53%  '

# Row 49
Set-TextCell ($ws.Cells.Item(49, 4)) '0.0506'
Set-TextCell ($ws.Cells.Item(49, 5)) '  -4.73%  '

# Row 50
Set-TextCell ($ws.Cells.Item(50, 4)) '19.33'
Set-TextCell ($ws.Cells.Item(50, 5)) '  -9.35%  '

# Row 51
Set-TextCell ($ws.Cells.Item(51, 4)) '0.565'
Set-TextCell ($ws.Cells.Item(51, 5)) '  -3.36%  '
